# Update countries & provincias Spain
#
# This script refreshes the Covid-19 "Pais" dashboard:
#   1) Bumps the "last updated" timestamp shown in A1.
#   2) Re-labels a handful of rows whose country changed place because the
#      source feed re-sorted/inserted entries (Libano/Mauritania/Albania,
#      Malaui/Republica de Africa Central, Aruba/Guyana).
#   3) Writes the newest totals (Casos totales, Nuevos casos, Casos activos,
#      Recuperados, Casos criticos, Muertes hoy, Muertes) for every country
#      row whose figures changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Refresh the "Datos actualizados..." banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 9 de Agosto de 2020 a las 19:53"

# 2) Fix up the country labels that shifted rows.
$labelUpdates = @(
    @{ Row = 98;  Name = "Libano" },
    @{ Row = 99;  Name = "Mauritania" },
    @{ Row = 100; Name = "Albania" },
    @{ Row = 108; Name = "Malaui" },
    @{ Row = 109; Name = "Republica de Africa Central" },
    @{ Row = 164; Name = "Aruba" },
    @{ Row = 165; Name = "Guyana" }
)

foreach ($item in $labelUpdates) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Name
}

# 3) Apply the refreshed numeric figures (columns B:H) for every changed row.
$rowData = @(
    @{ Row = 4;   B = 5167701; C = 17978; D = 2640081; E = 2362346; F = 0; G = 204; H = 165274 },
    @{ Row = 12;  B = 373056;  C = 2033;  D = 345826;  E = 17153;   F = 0; G = 66;  H = 10077 },
    @{ Row = 22;  B = 217167;  C = 271;   D = 197400;  E = 10506;   F = 0; G = 0;   H = 9261 },
    @{ Row = 31;  B = 94459;   C = 887;   D = 71605;   E = 16938;   F = 0; G = 0;   H = 5916 },
    @{ Row = 34;  B = 82670;   C = 346;   D = 57514;   E = 24556;   F = 0; G = 7;   H = 600 },
    @{ Row = 58;  B = 36603;   C = 152;   D = 32300;   E = 2317;    F = 0; G = 0;   H = 1986 },
    @{ Row = 59;  B = 35160;   C = 467;   D = 24506;   E = 9352;    F = 0; G = 9;   H = 1302 },
    @{ Row = 61;  B = 33237;   C = 1230;  D = 23347;   E = 9392;    F = 0; G = 18;  H = 498 },
    @{ Row = 65;  B = 26712;   C = 68;    D = 23364;   E = 1576;    F = 0; G = 0;   H = 1772 },
    @{ Row = 66;  B = 26436;   C = 599;   D = 12961;   E = 13055;   F = 0; G = 2;   H = 420 },
    @{ Row = 98;  B = 6517;    C = 294;   D = 2127;    E = 4312;    F = 0; G = 0;   H = 78 },
    @{ Row = 99;  B = 6510;    C = 0;     D = 5527;    E = 826;     F = 0; G = 0;   H = 157 },
    @{ Row = 100; B = 6411;    C = 136;   D = 3342;    E = 2870;    F = 0; G = 6;   H = 199 },
    @{ Row = 105; B = 5041;    C = 143;   D = 2804;    E = 2218;    F = 0; G = 0;   H = 19 },
    @{ Row = 108; B = 4658;    C = 34;    D = 2375;    E = 2137;    F = 0; G = 3;   H = 146 },
    @{ Row = 109; B = 4641;    C = 0;     D = 1716;    E = 2866;    F = 0; G = 0;   H = 59 },
    @{ Row = 125; B = 2470;    C = 7;     D = 1175;    E = 1248;    F = 0; G = 0;   H = 47 },
    @{ Row = 142; B = 1283;    C = 16;    D = 1115;    E = 161;     F = 0; G = 1;   H = 7 },
    @{ Row = 164; B = 563;     C = 54;    D = 114;     E = 446;     F = 0; G = 0;   H = 3 },
    @{ Row = 165; B = 554;     C = 0;     D = 189;     E = 343;     F = 0; G = 0;   H = 22 }
)

foreach ($item in $rowData) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
}
